$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.986.11'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.261.75'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.63'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("E6").Value = '  +2.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.29'
$ws.Range("E7").Value = '  +4.71%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.447'
$ws.Range("E9").Value = '  +6.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0997'
$ws.Range("E10").Value = '  +5.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.21'
$ws.Range("E11").Value = '  -1.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '27.13'
$ws.Range("E12").Value = '  +14.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.603.52'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.08'
$ws.Range("E16").Value = '  +5.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.836'
$ws.Range("E17").Value = '  +3.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.266.15'
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.971.63'
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("E20").Value = '  +7.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.66'
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.09'
$ws.Range("E22").Value = '  -2.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.75'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("E25").Value = '  -3.85%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("E26").Value = '  -2.51%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.09'
$ws.Range("E27").Value = '  +2.77%  '
$ws.Range("E28").Value = '  +21.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.20'
$ws.Range("E29").Value = '  +0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.140'
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.91'
$ws.Range("E31").Value = '  +2.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.39'
$ws.Range("E32").Value = '  -4.14%  '
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0707'
$ws.Range("E34").Value = '  +7.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.78'
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.90'
$ws.Range("E36").Value = '  -3.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.79'
$ws.Range("E37").Value = '  +5.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.52'
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.31'
$ws.Range("E39").Value = '  -3.54%  '
$ws.Range("E40").Value = '  +4.36%  '
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.000223'
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0979'
$ws.Range("E43").Value = '  -1.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.43'
$ws.Range("E44").Value = '  +5.10%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.21'
$ws.Range("E45").Value = '  -5.48%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.38'
$ws.Range("E46").Value = '  +8.61%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.21'
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.21'
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.43'
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("E50").Value = '  +5.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.443.91'
$ws.Range("E51").Value = '  -1.77%  '
